# 6.7 Add Investigate Panel
# Consolidate the misspelled "Quan-Regualr" / "Quna-Regular" avatar names
# into the correctly spelled "Quan-Regular", and update the active
# selection on Sheet1 to C17.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix the misspelled avatar names in column C so they all read "Quan-Regular".
$ws.Range("C3").Value = "Quan-Regular"
$ws.Range("C8").Value = "Quan-Regular"
$ws.Range("C9").Value = "Quan-Regular"
$ws.Range("C14").Value = "Quan-Regular"
$ws.Range("C16").Value = "Quan-Regular"

# Update the current selection to C17.
[void]$ws.Range("C17").Select()
